$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price observation was recorded. It belongs chronologically
# between the existing row 14 and the (old) row 15, so insert a fresh row
# at position 15 and push everything below it down by one.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new observation.
$ws.Cells.Item(15, 1).Value  = 5
$ws.Cells.Item(15, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(15, 3).Value  = "Maule"
$ws.Cells.Item(15, 4).Value  = 44802
$ws.Cells.Item(15, 5).Value  = 7
$ws.Cells.Item(15, 6).Value  = 100112043
$ws.Cells.Item(15, 7).Value  = "Pepino dulce"
$ws.Cells.Item(15, 8).Value  = "Cultivar IV Región"
$ws.Cells.Item(15, 9).Value  = "Primera"
$ws.Cells.Item(15, 10).Value = 500
$ws.Cells.Item(15, 11).Value = 14000
$ws.Cells.Item(15, 12).Value = 14000
$ws.Cells.Item(15, 13).Value = 14000
$ws.Cells.Item(15, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(15, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(15, 16).Value = 778
$ws.Cells.Item(15, 17).Value = 18
$ws.Cells.Item(15, 18).Value = "Hortaliza"
